# Data Sheet.xlsx edit:
#  - Food!B2: 6000 -> 9000 (cascades through the B5/B6/B7 formulas)
#  - Active sheet switches from "District" to "Food", with Food's selection
#    moving to D6 (District keeps its D13 selection, just loses focus)

$wb = $excel.ActiveWorkbook

$food = $wb.Worksheets.Item("Food")
$food.Range("B2").Value = 9000

$food.Activate()
$food.Range("D6").Select()
